$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to text format so numeric-looking strings
# (e.g. "1.99", "0.999", thousand-separated prices) are preserved verbatim
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '71.928.17'
$ws.Cells.Item(2, 5).Value = '  -0.06%  '
$ws.Cells.Item(3, 4).Value = '2.683.60'
$ws.Cells.Item(3, 5).Value = '  +1.92%  '
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
$ws.Cells.Item(5, 4).Value = '597.37'
$ws.Cells.Item(5, 5).Value = '  -1.90%  '
$ws.Cells.Item(6, 4).Value = '174.13'
$ws.Cells.Item(6, 5).Value = '  -3.07%  '
$ws.Cells.Item(7, 5).Value = '  -0.07%  '
$ws.Cells.Item(8, 4).Value = '0.523'
$ws.Cells.Item(8, 5).Value = '  -0.45%  '
$ws.Cells.Item(9, 4).Value = '2.683.04'
$ws.Cells.Item(9, 5).Value = '  +1.99%  '
$ws.Cells.Item(10, 4).Value = '0.166'
$ws.Cells.Item(10, 5).Value = '  -2.24%  '
$ws.Cells.Item(11, 5).Value = '  +1.96%  '
$ws.Cells.Item(12, 4).Value = '0.353'
$ws.Cells.Item(12, 5).Value = '  +1.39%  '
$ws.Cells.Item(13, 4).Value = '4.98'
$ws.Cells.Item(13, 5).Value = '  -0.97%  '
$ws.Cells.Item(14, 4).Value = '3.173.22'
$ws.Cells.Item(14, 5).Value = '  +2.89%  '
$ws.Cells.Item(15, 2).Value = 'WrappedBTC'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(15, 4).Value = '71.793.64'
$ws.Cells.Item(15, 5).Value = '  -0.13%  '
$ws.Cells.Item(16, 2).Value = 'ShibaInu'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(16, 4).Value = '0.0000184'
$ws.Cells.Item(16, 5).Value = '  -2.46%  '
$ws.Cells.Item(17, 4).Value = '26.14'
$ws.Cells.Item(17, 5).Value = '  -1.50%  '
$ws.Cells.Item(18, 4).Value = '2.689.23'
$ws.Cells.Item(18, 5).Value = '  +2.13%  '
$ws.Cells.Item(19, 4).Value = '12.21'
$ws.Cells.Item(19, 5).Value = '  +6.02%  '
$ws.Cells.Item(20, 4).Value = '8.08'
$ws.Cells.Item(20, 5).Value = '  +0.34%  '
$ws.Cells.Item(21, 4).Value = '370.64'
$ws.Cells.Item(21, 5).Value = '  -3.33%  '
$ws.Cells.Item(22, 4).Value = '4.17'
$ws.Cells.Item(22, 5).Value = '  +0.41%  '
$ws.Cells.Item(23, 4).Value = '1.99'
$ws.Cells.Item(23, 5).Value = '  -0.37%  '
$ws.Cells.Item(24, 4).Value = '72.10'
$ws.Cells.Item(24, 5).Value = '  -0.91%  '
$ws.Cells.Item(25, 5).Value = '  -0.10%  '
$ws.Cells.Item(26, 4).Value = '4.34'
$ws.Cells.Item(26, 5).Value = '  -2.63%  '
$ws.Cells.Item(27, 4).Value = '9.77'
$ws.Cells.Item(28, 4).Value = '2.819.67'
$ws.Cells.Item(28, 5).Value = '  +1.90%  '
$ws.Cells.Item(29, 5).Value = '  +0.38%  '
$ws.Cells.Item(30, 4).Value = '0.0₃0960'
$ws.Cells.Item(30, 5).Value = '  -0.88%  '
$ws.Cells.Item(31, 4).Value = '8.06'
$ws.Cells.Item(31, 5).Value = '  -0.10%  '
$ws.Cells.Item(32, 4).Value = '501.88'
$ws.Cells.Item(32, 5).Value = '  -8.10%  '
$ws.Cells.Item(33, 4).Value = '1.30'
$ws.Cells.Item(33, 5).Value = '  -2.59%  '
$ws.Cells.Item(34, 4).Value = '1.81'
$ws.Cells.Item(34, 5).Value = '  -1.15%  '
$ws.Cells.Item(35, 4).Value = '0.999'
$ws.Cells.Item(35, 5).Value = '  -0.12%  '
$ws.Cells.Item(36, 4).Value = '163.81'
$ws.Cells.Item(36, 5).Value = '  -1.40%  '
$ws.Cells.Item(37, 4).Value = '19.57'
$ws.Cells.Item(37, 5).Value = '  +1.61%  '
$ws.Cells.Item(38, 4).Value = '19.08'
$ws.Cells.Item(38, 5).Value = '  -0.21%  '
$ws.Cells.Item(39, 4).Value = '1.38'
$ws.Cells.Item(39, 5).Value = '  -1.47%  '
$ws.Cells.Item(40, 5).Value = '  -6.48%  '
$ws.Cells.Item(41, 5).Value = '  -4.22%  '
$ws.Cells.Item(42, 5).Value = '  -0.03%  '
$ws.Cells.Item(43, 4).Value = '5.00'
$ws.Cells.Item(43, 5).Value = '  -0.70%  '
$ws.Cells.Item(44, 4).Value = '0.333'
$ws.Cells.Item(44, 5).Value = '  +0.02%  '
$ws.Cells.Item(45, 4).Value = '2.54'
$ws.Cells.Item(45, 5).Value = '  -3.08%  '
$ws.Cells.Item(46, 4).Value = '156.57'
$ws.Cells.Item(46, 5).Value = '  +3.70%  '
$ws.Cells.Item(47, 4).Value = '39.33'
$ws.Cells.Item(47, 5).Value = '  -0.70%  '
$ws.Cells.Item(48, 4).Value = '0.559'
$ws.Cells.Item(48, 5).Value = '  +3.73%  '
$ws.Cells.Item(49, 4).Value = '3.73'
$ws.Cells.Item(49, 5).Value = '  +2.01%  '
$ws.Cells.Item(50, 4).Value = '1.75'
$ws.Cells.Item(50, 5).Value = '  +3.58%  '
$ws.Cells.Item(51, 4).Value = '0.0761'
$ws.Cells.Item(51, 5).Value = '  +0.86%  '
